$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.114.35"
$ws.Range("E2").Value = "  -0.64%  "
$ws.Range("D3").Value = "1.786.80"
$ws.Range("E3").Value = "  -2.67%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "223.83"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.46%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.550"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -1.19%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "32.71"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +2.50%  "
$ws.Range("E9").Value = "  -2.10%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0711"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.57%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0938"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.58%  "
$ws.Range("D12").Value = "2.045.46"
$ws.Range("E12").Value = "  -2.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.03"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +3.01%  "
$ws.Range("D14").Value = "1.789.34"
$ws.Range("E14").Value = "  -2.56%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.623"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -3.68%  "
$ws.Range("D16").Value = "34.066.98"
$ws.Range("E16").Value = "  -0.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.16"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -3.89%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.82"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -2.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.27"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.40%  "
$ws.Range("D20").Value = "0.0₃0789"
$ws.Range("E20").Value = "  -0.96%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.998"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  -4.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.09"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -3.87%  "
$ws.Range("E24").Value = "  -2.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "160.84"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "16.33"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -1.72%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.07"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.56%  "
$ws.Range("E28").Value = "  -2.38%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +0.01%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0513"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  -4.15%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.21"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("E32").Value = "  -3.58%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.51"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -2.22%  "
$ws.Range("E34").Value = "  -5.17%  "
$ws.Range("D35").Value = "1.396.03"
$ws.Range("E35").Value = "  -3.87%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.643"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -0.28%  "
$ws.Range("E37").Value = "  -1.61%  "
$ws.Range("E38").Value = "  -2.64%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.20"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +2.69%  "
$ws.Range("E40").Value = "  -0.36%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.915"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -5.39%  "
$ws.Range("E42").Value = "  -2.72%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "78.01"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -4.81%  "
$ws.Range("D44").Value = "0.0₆0143"
$ws.Range("E44").Value = "  +13.84%  "
$ws.Range("E45").Value = "  +2.19%  "
$ws.Range("B46").Value = "Kaspa"
$ws.Range("C46").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0498"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "108.13"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.47%  "
$ws.Range("B48").Value = "FraxShare"
$ws.Range("C48").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.87"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  -3.35%  "
$ws.Range("B49").Value = "InjectiveProtocol"
$ws.Range("C49").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "12.29"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.27%  "
$ws.Range("D50").Value = "1.945.04"
$ws.Range("E50").Value = "  -2.64%  "
$ws.Range("E51").Value = "  -0.06%  "
